$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# --- Update Hoja1!A1 text (conversion rates) ---
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.43 = 39778.49 pesos`n✅ 39778.49 pesos = 9.41 = 967.0 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update tasas sheet numeric cells ---
$wsTasas.Range("N10").Value = 106
$wsTasas.Range("O10").Value = 4216.52
$wsTasas.Range("N12").Value = 4229
$wsTasas.Range("O12").Value = 102.805
